$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.780.05'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.643.35'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.50%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '216.72'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  -0.09%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '19.14'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  -1.58%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0841'
$cell.Style = $origStyle
$ws.Range('D12').Value = '1.869.20'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '1.628.05'
$ws.Range('E13').Value = '  -0.88%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.16'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  -1.75%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.526'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.69%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '64.46'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -2.89%  '
$ws.Range('D17').Value = '26.791.44'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('E18').Value = '  -2.56%  '
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '214.17'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('E21').Value = '  -1.06%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.39'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +13.88%  '
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('E24').Value = '  -2.54%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '144.94'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('E27').Value = '  -2.09%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.05'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -0.81%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '15.68'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  +0.33%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.31'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').Value = '1.292.98'
$ws.Range('E34').Value = '  +0.26%  '
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.53'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  -1.41%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +1.17%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0174'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -6.04%  '
$ws.Range('E38').Value = '  +2.26%  '
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('E40').Value = '  +0.43%  '
$ws.Range('E41').Value = '  -0.29%  '
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.23'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  -0.41%  '
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').Value = '1.795.14'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '91.32'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -2.79%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '60.32'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('E47').Value = '  -1.19%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0520'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.68'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0980'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +0.58%  '
